# Applies the 2025-12-30 Betfair Back/Lay odds refresh:
#   - updates odds across several existing rows with newer quotes
#   - rows 18/20 (Altrincham v York City, Forest Green v Truro City) swap position
#   - rows 26/27 (Arsenal v Aston Villa, Man Utd v Wolves) swap position (with refreshed odds)
#   - row 28 (Torreense v Lusitania Futebol Clube) is removed, shrinking the sheet to A1:AO27
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Al-Adalh vs Al Jabalain
$ws.Cells.Item(2, 14).Value = 1.17  # N2
$ws.Cells.Item(2, 16).Value = 1.17  # P2
$ws.Cells.Item(2, 18).Value = 1.15  # R2

# Row 3: Al Ahli vs Al-Feiha
$ws.Cells.Item(3, 6).Value = 1.48  # F3
$ws.Cells.Item(3, 7).Value = 1.57  # G3
$ws.Cells.Item(3, 8).Value = 1.09  # H3
$ws.Cells.Item(3, 9).Value = 22  # I3
$ws.Cells.Item(3, 10).Value = 4.6  # J3
$ws.Cells.Item(3, 11).Value = 950  # K3
$ws.Cells.Item(3, 12).Value = 1.01  # L3
$ws.Cells.Item(3, 13).Value = 1.04  # M3
$ws.Cells.Item(3, 14).Value = 1.96  # N3
$ws.Cells.Item(3, 15).Value = 1.26  # O3
$ws.Cells.Item(3, 16).Value = 1.95  # P3
$ws.Cells.Item(3, 17).Value = 1.73  # Q3
$ws.Cells.Item(3, 18).Value = 1.18  # R3
$ws.Cells.Item(3, 19).Value = 1.73  # S3
$ws.Cells.Item(3, 20).Value = 1.89  # T3
$ws.Cells.Item(3, 21).Value = 1.79  # U3
$ws.Cells.Item(3, 22).Value = 1.12  # V3
$ws.Cells.Item(3, 23).Value = 2.74  # W3
$ws.Cells.Item(3, 24).Value = 1000  # X3
$ws.Cells.Item(3, 25).Value = 1000  # Y3
$ws.Cells.Item(3, 26).Value = 1000  # Z3
$ws.Cells.Item(3, 27).Value = 1000  # AA3
$ws.Cells.Item(3, 28).Value = 1000  # AB3
$ws.Cells.Item(3, 29).Value = 1000  # AC3
$ws.Cells.Item(3, 30).Value = 1000  # AD3
$ws.Cells.Item(3, 31).Value = 1000  # AE3
$ws.Cells.Item(3, 32).Value = 1000  # AF3
$ws.Cells.Item(3, 33).Value = 1000  # AG3
$ws.Cells.Item(3, 34).Value = 1000  # AH3
$ws.Cells.Item(3, 35).Value = 1000  # AI3
$ws.Cells.Item(3, 36).Value = 1000  # AJ3
$ws.Cells.Item(3, 37).Value = 1000  # AK3
$ws.Cells.Item(3, 38).Value = 1000  # AL3
$ws.Cells.Item(3, 39).Value = 1000  # AM3
$ws.Cells.Item(3, 40).Value = 1000  # AN3
$ws.Cells.Item(3, 41).Value = 1000  # AO3

# Row 4: Al-Ettifaq vs Al Nassr
$ws.Cells.Item(4, 7).Value = 1000  # G4
$ws.Cells.Item(4, 8).Value = 1.2  # H4
$ws.Cells.Item(4, 16).Value = 2.6  # P4
$ws.Cells.Item(4, 17).Value = 1.4  # Q4

# Row 6: Chelsea vs Bournemouth
$ws.Cells.Item(6, 8).Value = 6  # H6
$ws.Cells.Item(6, 20).Value = 1.67  # T6
$ws.Cells.Item(6, 24).Value = 25  # X6
$ws.Cells.Item(6, 25).Value = 29  # Y6
$ws.Cells.Item(6, 27).Value = 170  # AA6
$ws.Cells.Item(6, 32).Value = 11.5  # AF6
$ws.Cells.Item(6, 37).Value = 15  # AK6
$ws.Cells.Item(6, 38).Value = 28  # AL6

# Row 7: West Ham vs Brighton
$ws.Cells.Item(7, 8).Value = 2.24  # H7
$ws.Cells.Item(7, 9).Value = 2.26  # I7
$ws.Cells.Item(7, 10).Value = 3.75  # J7
$ws.Cells.Item(7, 24).Value = 18.5  # X7
$ws.Cells.Item(7, 35).Value = 36  # AI7
$ws.Cells.Item(7, 39).Value = 1000  # AM7

# Row 8: Burnley vs Newcastle
$ws.Cells.Item(8, 9).Value = 1.72  # I8
$ws.Cells.Item(8, 28).Value = 19.5  # AB8
$ws.Cells.Item(8, 33).Value = 21  # AG8

# Row 9: Nottm Forest vs Everton
$ws.Cells.Item(9, 15).Value = 1.43  # O9
$ws.Cells.Item(9, 16).Value = 1.74  # P9
$ws.Cells.Item(9, 17).Value = 2.3  # Q9
$ws.Cells.Item(9, 41).Value = 80  # AO9

# Row 13: Hibernian vs Aberdeen
$ws.Cells.Item(13, 6).Value = 1.9  # F13
$ws.Cells.Item(13, 7).Value = 1.93  # G13
$ws.Cells.Item(13, 8).Value = 4.2  # H13
$ws.Cells.Item(13, 10).Value = 3.9  # J13
$ws.Cells.Item(13, 11).Value = 4.1  # K13

# Row 14: Scunthorpe vs Gateshead
$ws.Cells.Item(14, 17).Value = 1.48  # Q14

# Row 17: Tamworth FC vs FC Halifax Town
$ws.Cells.Item(17, 6).Value = 2.88  # F17
$ws.Cells.Item(17, 8).Value = 2.42  # H17
$ws.Cells.Item(17, 9).Value = 2.58  # I17
$ws.Cells.Item(17, 10).Value = 3.6  # J17
$ws.Cells.Item(17, 16).Value = 1.94  # P17

# Row 18: Forest Green vs Truro City
$ws.Cells.Item(18, 4).Value = "Forest Green"  # D18
$ws.Cells.Item(18, 5).Value = "Truro City"  # E18
$ws.Cells.Item(18, 6).Value = 1.21  # F18
$ws.Cells.Item(18, 7).Value = 1000  # G18
$ws.Cells.Item(18, 8).Value = 1.04  # H18
$ws.Cells.Item(18, 9).Value = 1000  # I18
$ws.Cells.Item(18, 10).Value = 5.1  # J18
$ws.Cells.Item(18, 11).Value = 950  # K18
$ws.Cells.Item(18, 16).Value = 2.22  # P18
$ws.Cells.Item(18, 17).Value = 1.47  # Q18

# Row 19: Rochdale vs Hartlepool
$ws.Cells.Item(19, 6).Value = 1.58  # F19

# Row 20: Altrincham vs York City
$ws.Cells.Item(20, 4).Value = "Altrincham"  # D20
$ws.Cells.Item(20, 5).Value = "York City"  # E20
$ws.Cells.Item(20, 6).Value = 5.3  # F20
$ws.Cells.Item(20, 7).Value = 7.4  # G20
$ws.Cells.Item(20, 8).Value = 1.54  # H20
$ws.Cells.Item(20, 9).Value = 1.66  # I20
$ws.Cells.Item(20, 10).Value = 4.6  # J20
$ws.Cells.Item(20, 11).Value = 5.5  # K20
$ws.Cells.Item(20, 16).Value = 2.5  # P20
$ws.Cells.Item(20, 17).Value = 1.01  # Q20

# Row 21: Yeovil vs Eastleigh
$ws.Cells.Item(21, 6).Value = 2.52  # F21
$ws.Cells.Item(21, 8).Value = 2.86  # H21
$ws.Cells.Item(21, 9).Value = 3.55  # I21

# Row 22: Boreham Wood vs Solihull Moors
$ws.Cells.Item(22, 7).Value = 1.63  # G22
$ws.Cells.Item(22, 10).Value = 4.5  # J22

# Row 23: Aldershot vs Wealdstone
$ws.Cells.Item(23, 6).Value = 2.26  # F23
$ws.Cells.Item(23, 8).Value = 2.84  # H23

# Row 26: Man Utd vs Wolves
$ws.Cells.Item(26, 4).Value = "Man Utd"  # D26
$ws.Cells.Item(26, 5).Value = "Wolves"  # E26
$ws.Cells.Item(26, 6).Value = 1.38  # F26
$ws.Cells.Item(26, 7).Value = 1.39  # G26
$ws.Cells.Item(26, 8).Value = 9.6  # H26
$ws.Cells.Item(26, 9).Value = 10  # I26
$ws.Cells.Item(26, 10).Value = 5.7  # J26
$ws.Cells.Item(26, 11).Value = 5.8  # K26
$ws.Cells.Item(26, 13).Value = 1.04  # M26
$ws.Cells.Item(26, 14).Value = 0  # N26
$ws.Cells.Item(26, 15).Value = 0  # O26
$ws.Cells.Item(26, 16).Value = 2.52  # P26
$ws.Cells.Item(26, 17).Value = 1.62  # Q26
$ws.Cells.Item(26, 18).Value = 0  # R26
$ws.Cells.Item(26, 19).Value = 0  # S26
$ws.Cells.Item(26, 20).Value = 1.94  # T26
$ws.Cells.Item(26, 21).Value = 2  # U26
$ws.Cells.Item(26, 24).Value = 26  # X26
$ws.Cells.Item(26, 25).Value = 34  # Y26
$ws.Cells.Item(26, 26).Value = 90  # Z26
$ws.Cells.Item(26, 27).Value = 370  # AA26
$ws.Cells.Item(26, 28).Value = 10  # AB26
$ws.Cells.Item(26, 29).Value = 13  # AC26
$ws.Cells.Item(26, 30).Value = 36  # AD26
$ws.Cells.Item(26, 31).Value = 170  # AE26
$ws.Cells.Item(26, 32).Value = 8.8  # AF26
$ws.Cells.Item(26, 33).Value = 10.5  # AG26
$ws.Cells.Item(26, 34).Value = 26  # AH26
$ws.Cells.Item(26, 35).Value = 110  # AI26
$ws.Cells.Item(26, 36).Value = 11.5  # AJ26
$ws.Cells.Item(26, 37).Value = 14  # AK26
$ws.Cells.Item(26, 38).Value = 34  # AL26
$ws.Cells.Item(26, 39).Value = 140  # AM26
$ws.Cells.Item(26, 40).Value = 5.2  # AN26
$ws.Cells.Item(26, 41).Value = 1000  # AO26

# Row 27: Arsenal vs Aston Villa
$ws.Cells.Item(27, 4).Value = "Arsenal"  # D27
$ws.Cells.Item(27, 5).Value = "Aston Villa"  # E27
$ws.Cells.Item(27, 6).Value = 1.5  # F27
$ws.Cells.Item(27, 7).Value = 1.51  # G27
$ws.Cells.Item(27, 8).Value = 8  # H27
$ws.Cells.Item(27, 9).Value = 8.2  # I27
$ws.Cells.Item(27, 10).Value = 4.6  # J27
$ws.Cells.Item(27, 11).Value = 4.7  # K27
$ws.Cells.Item(27, 13).Value = 1.07  # M27
$ws.Cells.Item(27, 14).Value = 3.95  # N27
$ws.Cells.Item(27, 15).Value = 1.32  # O27
$ws.Cells.Item(27, 16).Value = 2  # P27
$ws.Cells.Item(27, 17).Value = 1.96  # Q27
$ws.Cells.Item(27, 18).Value = 1.37  # R27
$ws.Cells.Item(27, 19).Value = 3.5  # S27
$ws.Cells.Item(27, 20).Value = 2.18  # T27
$ws.Cells.Item(27, 21).Value = 1.8  # U27
$ws.Cells.Item(27, 24).Value = 15  # X27
$ws.Cells.Item(27, 25).Value = 22  # Y27
$ws.Cells.Item(27, 26).Value = 80  # Z27
$ws.Cells.Item(27, 27).Value = 500  # AA27
$ws.Cells.Item(27, 28).Value = 8  # AB27
$ws.Cells.Item(27, 29).Value = 10  # AC27
$ws.Cells.Item(27, 30).Value = 29  # AD27
$ws.Cells.Item(27, 31).Value = 160  # AE27
$ws.Cells.Item(27, 32).Value = 8.2  # AF27
$ws.Cells.Item(27, 33).Value = 10  # AG27
$ws.Cells.Item(27, 34).Value = 28  # AH27
$ws.Cells.Item(27, 35).Value = 140  # AI27
$ws.Cells.Item(27, 36).Value = 12.5  # AJ27
$ws.Cells.Item(27, 37).Value = 16.5  # AK27
$ws.Cells.Item(27, 38).Value = 44  # AL27
$ws.Cells.Item(27, 39).Value = 170  # AM27
$ws.Cells.Item(27, 40).Value = 8.2  # AN27
$ws.Cells.Item(27, 41).Value = 260  # AO27

# Row 28 (Torreense vs Lusitania Futebol Clube) no longer present in the refreshed feed;
# delete it so the sheet shrinks back to A1:AO27.
$ws.Rows.Item(28).Delete()
